$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.951738
$ws.Range("H2").Value = 5.855214
$ws.Range("I2").Value = 0.1200622759026226
$ws.Range("J2").Value = 0.1200622759026226
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.632955
$ws.Range("N2").Value = 94.898865
$ws.Range("O2").Value = 0.5000340016649593
$ws.Range("P2").Value = 0.5000340016649593
$ws.Range("Q2").Value = 61.73924032579
$ws.Range("R2").Value = 555.65316293211
$ws.Range("S2").Value = 0.06003522026859078
$ws.Range("T2").Value = 0.06003522026859078

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.951738
$ws.Range("H3").Value = 5.855214
$ws.Range("I3").Value = 0.1200622759026226
$ws.Range("J3").Value = 0.1200622759026226
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.622273333333332
$ws.Range("N3").Value = 25.86682
$ws.Range("O3").Value = 0.1362955132808722
$ws.Range("P3").Value = 0.1362955132808722
$ws.Range("Q3").Value = 16.82841851105333
$ws.Range("R3").Value = 151.45576659948
$ws.Range("S3").Value = 0.01636394951981764
$ws.Range("T3").Value = 0.01636394951981764

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.951738
$ws.Range("H4").Value = 5.855214
$ws.Range("I4").Value = 0.1200622759026226
$ws.Range("J4").Value = 0.1200622759026226
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.70876033333333
$ws.Range("N4").Value = 56.12628100000001
$ws.Range("O4").Value = 0.2957364019791172
$ws.Range("P4").Value = 0.2957364019791172
$ws.Range("Q4").Value = 36.51459847545934
$ws.Range("R4").Value = 328.6313862791341
$ws.Range("S4").Value = 0.03550678548886566
$ws.Range("T4").Value = 0.03550678548886567

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.951738
$ws.Range("H5").Value = 5.855214
$ws.Range("I5").Value = 0.1200622759026226
$ws.Range("J5").Value = 0.1200622759026226
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.297619333333333
$ws.Range("N5").Value = 12.892858
$ws.Range("O5").Value = 0.06793408307505136
$ws.Range("P5").Value = 0.06793408307505136
$ws.Range("Q5").Value = 8.387826962401334
$ws.Range("R5").Value = 75.49044266161201
$ws.Range("S5").Value = 0.0081563206253485
$ws.Range("T5").Value = 0.0081563206253485

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 9.654910333333332
$ws.Range("H6").Value = 28.964731
$ws.Range("I6").Value = 0.5939273141455197
$ws.Range("J6").Value = 0.5939273141455197
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.632955
$ws.Range("N6").Value = 94.898865
$ws.Range("O6").Value = 0.5000340016649593
$ws.Range("P6").Value = 0.5000340016649593
$ws.Range("Q6").Value = 305.4133441033683
$ws.Range("R6").Value = 2748.720096930315
$ws.Range("S6").Value = 0.2969838515903055
$ws.Range("T6").Value = 0.2969838515903055

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 9.654910333333332
$ws.Range("H7").Value = 28.964731
$ws.Range("I7").Value = 0.5939273141455197
$ws.Range("J7").Value = 0.5939273141455197
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.622273333333332
$ws.Range("N7").Value = 25.86682
$ws.Range("O7").Value = 0.1362955132808722
$ws.Range("P7").Value = 0.1362955132808722
$ws.Range("Q7").Value = 83.24727590282443
$ws.Range("R7").Value = 749.2254831254198
$ws.Range("S7").Value = 0.08094962813299345
$ws.Range("T7").Value = 0.08094962813299345

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.654910333333332
$ws.Range("H8").Value = 28.964731
$ws.Range("I8").Value = 0.5939273141455197
$ws.Range("J8").Value = 0.5939273141455197
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 18.70876033333333
$ws.Range("N8").Value = 56.12628100000001
$ws.Range("O8").Value = 0.2957364019791172
$ws.Range("P8").Value = 0.2957364019791172
$ws.Range("Q8").Value = 180.6314034661568
$ws.Range("R8").Value = 1625.682631195411
$ws.Range("S8").Value = 0.1756459269225168
$ws.Range("T8").Value = 0.1756459269225168

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.654910333333332
$ws.Range("H9").Value = 28.964731
$ws.Range("I9").Value = 0.5939273141455197
$ws.Range("J9").Value = 0.5939273141455197
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.297619333333333
$ws.Range("N9").Value = 12.892858
$ws.Range("O9").Value = 0.06793408307505136
$ws.Range("P9").Value = 0.06793408307505136
$ws.Range("Q9").Value = 41.49312931013311
$ws.Range("R9").Value = 373.438163791198
$ws.Range("S9").Value = 0.04034790749970386
$ws.Range("T9").Value = 0.04034790749970386

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.050938666666667
$ws.Range("H10").Value = 3.152816
$ws.Range("I10").Value = 0.06464909129917419
$ws.Range("J10").Value = 0.06464909129917419
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.632955
$ws.Range("N10").Value = 94.898865
$ws.Range("O10").Value = 0.5000340016649593
$ws.Range("P10").Value = 0.5000340016649593
$ws.Range("Q10").Value = 33.24429555042666
$ws.Range("R10").Value = 299.19865995384
$ws.Range("S10").Value = 0.03232674382632937
$ws.Range("T10").Value = 0.03232674382632937

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.050938666666667
$ws.Range("H11").Value = 3.152816
$ws.Range("I11").Value = 0.06464909129917419
$ws.Range("J11").Value = 0.06464909129917419
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.622273333333332
$ws.Range("N11").Value = 25.86682
$ws.Range("O11").Value = 0.1362955132808722
$ws.Range("P11").Value = 0.1362955132808722
$ws.Range("Q11").Value = 9.061480440568888
$ws.Range("R11").Value = 81.55332396511999
$ws.Range("S11").Value = 0.008811381081762916
$ws.Range("T11").Value = 0.008811381081762916

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.050938666666667
$ws.Range("H12").Value = 3.152816
$ws.Range("I12").Value = 0.06464909129917419
$ws.Range("J12").Value = 0.06464909129917419
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.70876033333333
$ws.Range("N12").Value = 56.12628100000001
$ws.Range("O12").Value = 0.2957364019791172
$ws.Range("P12").Value = 0.2957364019791172
$ws.Range("Q12").Value = 19.66175963969956
$ws.Range("R12").Value = 176.955836757296
$ws.Range("S12").Value = 0.01911908965203723
$ws.Range("T12").Value = 0.01911908965203723

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.050938666666667
$ws.Range("H13").Value = 3.152816
$ws.Range("I13").Value = 0.06464909129917419
$ws.Range("J13").Value = 0.06464909129917419
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.297619333333333
$ws.Range("N13").Value = 12.892858
$ws.Range("O13").Value = 0.06793408307505136
$ws.Range("P13").Value = 0.06793408307505136
$ws.Range("Q13").Value = 4.516534332014222
$ws.Range("R13").Value = 40.648808988128
$ws.Range("S13").Value = 0.00439187673904468
$ws.Range("T13").Value = 0.00439187673904468

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.59846
$ws.Range("H14").Value = 10.79538
$ws.Range("I14").Value = 0.2213613186526835
$ws.Range("J14").Value = 0.2213613186526835
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.632955
$ws.Range("N14").Value = 94.898865
$ws.Range("O14").Value = 0.5000340016649593
$ws.Range("P14").Value = 0.5000340016649593
$ws.Range("Q14").Value = 113.8299232493
$ws.Range("R14").Value = 1024.4693092437
$ws.Range("S14").Value = 0.1106881859797335
$ws.Range("T14").Value = 0.1106881859797335

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.59846
$ws.Range("H15").Value = 10.79538
$ws.Range("I15").Value = 0.2213613186526835
$ws.Range("J15").Value = 0.2213613186526835
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 8.622273333333332
$ws.Range("N15").Value = 25.86682
$ws.Range("O15").Value = 0.1362955132808722
$ws.Range("P15").Value = 0.1362955132808722
$ws.Range("Q15").Value = 31.02690569906666
$ws.Range("R15").Value = 279.2421512916
$ws.Range("S15").Value = 0.03017055454629821
$ws.Range("T15").Value = 0.03017055454629821

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.59846
$ws.Range("H16").Value = 10.79538
$ws.Range("I16").Value = 0.2213613186526835
$ws.Range("J16").Value = 0.2213613186526835
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 18.70876033333333
$ws.Range("N16").Value = 56.12628100000001
$ws.Range("O16").Value = 0.2957364019791172
$ws.Range("P16").Value = 0.2957364019791172
$ws.Range("Q16").Value = 67.32272570908667
$ws.Range("R16").Value = 605.9045313817801
$ws.Range("S16").Value = 0.06546459991569746
$ws.Range("T16").Value = 0.06546459991569747

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.59846
$ws.Range("H17").Value = 10.79538
$ws.Range("I17").Value = 0.2213613186526835
$ws.Range("J17").Value = 0.2213613186526835
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.297619333333333
$ws.Range("N17").Value = 12.892858
$ws.Range("O17").Value = 0.06793408307505136
$ws.Range("P17").Value = 0.06793408307505136
$ws.Range("Q17").Value = 15.46481126622667
$ws.Range("R17").Value = 139.18330139604
$ws.Range("S17").Value = 0.01503797821095432
$ws.Range("T17").Value = 0.01503797821095432
